# Juno: check in to OLPRODLOC.
# Translate the "Charger sales report" worksheet from English to French.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet (tab name) to the French title.
$ws.Name = "Rapport de ventes"

# --- Translate header row (row 1) ---
$ws.Range("A1").Value = "Année/trimestre"
$ws.Range("B1").Value = "Midwest"
$ws.Range("C1").Value = "Montagne"
$ws.Range("D1").Value = "Northeast"
$ws.Range("E1").Value = "Sud"
$ws.Range("F1").Value = "Sud-est"
$ws.Range("G1").Value = "Ouest"

# --- Translate the Year-Quarter labels in column A (rows 2-9) ---
$ws.Range("A2").Value = "2022-T1"
$ws.Range("A3").Value = "2022-T2"
$ws.Range("A4").Value = "2022-T3"
$ws.Range("A5").Value = "2022-T4"
$ws.Range("A6").Value = "2023-T1"
$ws.Range("A7").Value = "2023-T2"
$ws.Range("A8").Value = "2023-T3"
$ws.Range("A9").Value = "2023-T4"

# --- Re-apply consistent font formatting (Aptos Narrow, 11pt, black) across
#     the whole used range of text labels, as rich text runs. Applying the
#     font in two passes (rather than the full run in one call) keeps the
#     run merged as a single formatted span per cell. ---
$labelCells = @("A1","B1","C1","D1","E1","F1","G1","A2","A3","A4","A5","A6","A7","A8","A9")
foreach ($addr in $labelCells) {
    $cell = $ws.Range($addr)
    $len = $cell.Characters().Count
    $split = [Math]::Max(1, $len - 1)

    $cell.Characters(1, $split).Font.Name = "Aptos Narrow"
    $cell.Characters(1, $split).Font.Size = 11
    $cell.Characters(1, $split).Font.ColorIndex = 1

    if ($len -gt $split) {
        $cell.Characters($split + 1, $len - $split).Font.Name = "Aptos Narrow"
        $cell.Characters($split + 1, $len - $split).Font.Size = 11
        $cell.Characters($split + 1, $len - $split).Font.ColorIndex = 1
    }
}
